$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: simple in-place updates to nombre_aides (C) and montant_total (D) ---
# (rows that do not shift position)
$simpleUpdates = @(
    @(2, '187', '442016.00'),
    @(3, '1000', '3179764.33'),
    @(4, '413', '1675698.25'),
    @(5, '115', '547128.09'),
    @(21, '52', '126800.00'),
    @(22, '326', '953679.20'),
    @(24, '41', '207233.00'),
    @(25, '8', '67000.00'),
    @(26, '19', '47500.00'),
    @(33, '106', '295673.00'),
    @(34, '560', '1823794.47'),
    @(35, '223', '1118288.11'),
    @(51, '583', '2008915.52'),
    @(55, '22', '66220.65'),
    @(63, '15', '34861.00'),
    @(64, '71', '179849.69'),
    @(65, '33', '110535.00'),
    @(66, '11', '54027.00')
)
foreach ($u in $simpleUpdates) {
    $r = $u[0]
    $ws.Cells.Item($r, 3).Value = "'" + $u[1]
    $ws.Cells.Item($r, 4).Value = "'" + $u[2]
}

# --- Step 2: insert two new blank rows for the new "La Réunion" classes ---
$ws.Range("67:68").Insert()

# --- Step 3: populate the two newly inserted rows ---
$newRows = @(
    @(67, 'Fonds de solidarité', 'VOLET2', '3', '10500.00', '04', 'La Réunion', '11', '10 à 19 salariés'),
    @(68, 'Fonds de solidarité', 'VOLET2', '4', '9500.00', '04', 'La Réunion', 'NN', 'Etablissement non employeur')
)
foreach ($nr in $newRows) {
    $r = $nr[0]
    $ws.Cells.Item($r, 1).Value = $nr[1]
    $ws.Cells.Item($r, 2).Value = $nr[2]
    $ws.Cells.Item($r, 3).Value = "'" + $nr[3]
    $ws.Cells.Item($r, 4).Value = "'" + $nr[4]
    $ws.Cells.Item($r, 5).Value = "'" + $nr[5]
    $ws.Cells.Item($r, 6).Value = $nr[6]
    $ws.Cells.Item($r, 7).Value = "'" + $nr[7]
    $ws.Cells.Item($r, 8).Value = $nr[8]
}

# --- Step 4: updates to rows that shifted down by 2 (Nouvelle-Aquitaine, PACA blocks) ---
$postShiftUpdates = @(
    @(81, '228', '581326.09'),
    @(82, '885', '2826547.26'),
    @(83, '336', '1357303.79'),
    @(84, '117', '571984.52'),
    @(85, '29', '174080.04'),
    @(86, '34', '74500.00'),
    @(99, '295', '774979.43'),
    @(100, '1212', '3691343.89'),
    @(101, '451', '1846594.02'),
    @(102, '121', '563996.00'),
    @(103, '34', '219157.00')
)
foreach ($u in $postShiftUpdates) {
    $r = $u[0]
    $ws.Cells.Item($r, 3).Value = "'" + $u[1]
    $ws.Cells.Item($r, 4).Value = "'" + $u[2]
}
